$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '70.840.93'
$ws.Range("E2").Value = '  -1.57%  '
$ws.Range("D3").Value = '2.571.77'
$ws.Range("E3").Value = '  -4.24%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.94'
$ws.Range("E5").Value = '  -3.12%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.26'
$ws.Range("E6").Value = '  -1.68%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  -2.21%  '
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("D10").Value = '2.571.24'
$ws.Range("E10").Value = '  -4.18%  '
$ws.Range("E11").Value = '  -0.14%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.350'
$ws.Range("E12").Value = '  -1.14%  '
$ws.Range("E13").Value = '  -2.69%  '
$ws.Range("E14").Value = '  +0.38%  '
$ws.Range("D15").Value = '3.041.88'
$ws.Range("E15").Value = '  -4.29%  '
$ws.Range("D16").Value = '70.717.54'
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '25.29'
$ws.Range("E17").Value = '  -3.43%  '
$ws.Range("D18").Value = '2.568.38'
$ws.Range("E18").Value = '  -4.29%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.84'
$ws.Range("E19").Value = '  -2.96%  '
$ws.Range("B20").Value = 'BitcoinCash'
$ws.Range("C20").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '364.64'
$ws.Range("E20").Value = '  -1.90%  '
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.50'
$ws.Range("E21").Value = '  -7.56%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.00'
$ws.Range("E22").Value = '  -4.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '2.02'
$ws.Range("E23").Value = '  +1.25%  '
$ws.Range("E24").Value = '  +0.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.28'
$ws.Range("E25").Value = '  -2.78%  '
$ws.Range("E26").Value = '  -4.38%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.37'
$ws.Range("E27").Value = '  -4.22%  '
$ws.Range("E28").Value = '  -4.24%  '
$ws.Range("E29").Value = '  +0.26%  '
$ws.Range("E30").Value = '  -2.66%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.86'
$ws.Range("E31").Value = '  -2.55%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '487.36'
$ws.Range("E32").Value = '  -2.33%  '
$ws.Range("E33").Value = '  +0.96%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.77'
$ws.Range("E34").Value = '  -2.36%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  -0.05%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '157.22'
$ws.Range("E36").Value = '  -4.08%  '
$ws.Range("E37").Value = '  +6.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '18.82'
$ws.Range("E38").Value = '  -3.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.85'
$ws.Range("E39").Value = '  -1.37%  '
$ws.Range("E40").Value = '  -2.92%  '
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("E42").Value = '  -4.93%  '
$ws.Range("E43").Value = '  -1.92%  '
$ws.Range("E44").Value = '  -4.19%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.322'
$ws.Range("E45").Value = '  -3.44%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '38.56'
$ws.Range("E46").Value = '  -2.10%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '146.50'
$ws.Range("E47").Value = '  -6.95%  '
$ws.Range("E48").Value = '  -3.87%  '
$ws.Range("E49").Value = '  -5.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.65'
$ws.Range("E50").Value = '  -5.94%  '
$ws.Range("E51").Value = '  -0.94%  '
